# Standardise "cost_variable" -> "cost_variable_om" in column C (rows 10-39)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("C10:C39")
$range.Value = "cost_variable_om"

# Update selection to match the edited range, as recorded in the workbook view
$ws.Range("C10:C39").Select() | Out-Null
